$d = $word.ActiveDocument

# Replace teacher name occurrences (both instances in the document)
$d.Content.Find.Execute("M. Umar Hadi, S.Pd.", $true, $false, $false, $false, $false, $true, 1, $false, "Chusnul Muawanah, S.T., MM.", 2)

# Replace NIP number occurrences (both instances in the document)
$d.Content.Find.Execute("19700417 198903 1003", $true, $false, $false, $false, $false, $true, 1, $false, "19800104 200901 2004", 2)
